$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.411322355270386
$ws.Range("B1").Value = 2.902334928512573
$ws.Range("C1").Value = 1.77352511882782
$ws.Range("D1").Value = 1.229152083396912
$ws.Range("E1").Value = 1.018756985664368
